# Daily attendance processing - 2025-10-15 22:20:06
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (column G) wherever the value is exactly
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Range("G$row")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
